$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 55, shifting rows 55:68 down to 56:69
$ws.Rows.Item(55).Insert()

# Fill the new row 55 with the new weekly record
$ws.Cells.Item(55, 1).Value = 11
$ws.Cells.Item(55, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(55, 3).Value = "Bíobío"
$ws.Cells.Item(55, 4).Value = 44551
$ws.Cells.Item(55, 5).Value = 8
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100101
$ws.Cells.Item(55, 8).Value = "Berries"
$ws.Cells.Item(55, 9).Value = 100101001
$ws.Cells.Item(55, 10).Value = "Arándano (blue)"
$ws.Cells.Item(55, 11).Value = "Sin especificar"
$ws.Cells.Item(55, 12).Value = "Primera"
$ws.Cells.Item(55, 13).Value = 250
$ws.Cells.Item(55, 14).Value = 2500
$ws.Cells.Item(55, 15).Value = 3000
$ws.Cells.Item(55, 16).Value = 2700
$ws.Cells.Item(55, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(55, 18).Value = "Región de Ñuble"
$ws.Cells.Item(55, 19).Value = 1350
$ws.Cells.Item(55, 20).Value = 2
